$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Updated solid_mg_mwh (B), liq_mg_mwh (C), gas_mg_mwh (D) values per eGRID subregion row.
# Row 5 - AZNM
$ws.Range("B5").Value = 21.474557569998417
$ws.Range("C5").Value = 0.0011139371050135263
$ws.Range("D5").Value = 0.033862719637723124

# Row 7 - FRCC
$ws.Range("B7").Value = 143.37057952933651
$ws.Range("C7").Value = 0.0077543764636170166
$ws.Range("D7").Value = 4.0946274621723573

# Row 8 - HIMS
$ws.Range("B8").Value = 73.146974683871193
$ws.Range("C8").Value = 0.0038028852378157933
$ws.Range("D8").Value = 1.1548169749744757

# Row 11 - MROE
$ws.Range("B11").Value = 26.379459841225277
$ws.Range("C11").Value = 0.0012907434249125433
$ws.Range("D11").Value = 0.56614439437360864

# Row 12 - MROW
$ws.Range("B12").Value = 242.02123048820144
$ws.Range("C12").Value = 0.021231579893888516
$ws.Range("D12").Value = 3.1185722888937435

# Row 14 - NWPP
$ws.Range("B14").Value = 161.21288435241371
$ws.Range("C14").Value = 0.0045729101723336373
$ws.Range("D14").Value = 0.97331036621811218

# Row 17 - NYUP
$ws.Range("B17").Value = 15.380522610799117
$ws.Range("C17").Value = 0.00019072841271542879
$ws.Range("D17").Value = 0.098111181622800803

# Row 18 - RFCE
$ws.Range("B18").Value = 142.06463302404146
$ws.Range("C18").Value = 0.0028652528870556069
$ws.Range("D18").Value = 4.6013404844227663

# Row 19 - RFCM
$ws.Range("B19").Value = 102.59925455126825
$ws.Range("C19").Value = 0.0034177471067292839
$ws.Range("D19").Value = 3.7123626844809352

# Row 20 - RFCW
$ws.Range("B20").Value = 749.15912282950364
$ws.Range("C20").Value = 0.058255920862844736
$ws.Range("D20").Value = 10.025801934582807

# Row 21 - RMPA
$ws.Range("B21").Value = 144.1690614397111
$ws.Range("C21").Value = 0.0075804305098319223
$ws.Range("D21").Value = 0.15834694711400502

# Row 22 - SPNO
$ws.Range("B22").Value = 36.462789248376531
$ws.Range("C22").Value = 0.00040007976612425747
$ws.Range("D22").Value = 0.55965815643999217

# Row 23 - SPSO
$ws.Range("B23").Value = 121.61017387057332
$ws.Range("C23").Value = 0.029927953734358855
$ws.Range("D23").Value = 3.2429363691050894

# Row 24 - SRMV
$ws.Range("B24").Value = 45.928649848411027
$ws.Range("C24").Value = 0.0030200931189827709
$ws.Range("D24").Value = 0.087749705103078959

# Row 25 - SRMW
$ws.Range("B25").Value = 74.721363344339252
$ws.Range("C25").Value = 0.00083189264424682059
$ws.Range("D25").Value = 3.9416560952647228

# Row 26 - SRSO
$ws.Range("B26").Value = 323.99435838102909
$ws.Range("C26").Value = 0.034819671609657897
$ws.Range("D26").Value = 20.86122161132851

# Row 27 - SRTV
$ws.Range("B27").Value = 361.5578598943535
$ws.Range("C27").Value = 0.01823919178431959
$ws.Range("D27").Value = 1.8834472767064019

# Row 28 - SRVC
$ws.Range("B28").Value = 192.15111546133838
$ws.Range("C28").Value = 0.0067217245421593332
$ws.Range("D28").Value = 0.41427121391883531
